$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Alt1")
$ws.Activate()

# Enter the (experimental / ultimately abandoned) dynamic-array formula that
# tries to "carry forward" the last valid date using SCAN + LAMBDA + ISNUMBER(DATEVALUE(...)).
# Anchored at L42, it spills down through L66 because the named range _na
# ('Alt1'!$C$3:$C$27) has 25 rows.
$ws.Range("L42").Formula = "=SCAN(,_na,LAMBDA(a,v,IF(ISNUMBER(DATEVALUE(TEXT(v,""dd-mmm-yyyy""))),v,a)))"

# Match the saved view state: scrolled down so row 40 is at the top, with
# G43 as the active selected cell.
$win = $excel.ActiveWindow
$win.ScrollRow = 40
$win.ScrollColumn = 1
$ws.Range("G43").Select()
